# akc rank checker done
#
# Sheet1: append a new "EA Name" entry ("BCSS Course") below the existing
# "Food Safety Course" row.
#
# AKC Rankings: the rank-checker ran again on 2025-11-04. Four more "Not
# Found" rows landed for "Food Safety Course" (rows 4-7), then three more
# runs finally resolved it at rank 7 (rows 8-10), and a first run for the
# newly added "BCSS Course" landed at rank 2 (row 11). The two older rows
# (2-3, still "Not Found") get their date column refreshed to 2025-11-04
# too. The bold/centered/bordered header style that used to be on row 1 is
# cleared.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: add "BCSS Course" as a new EA Name entry.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A3").Value = "BCSS Course"
$ws1.Range("A4").Select() | Out-Null

# ---------------------------------------------------------------------
# AKC Rankings: refresh dates, drop the header styling, append new runs.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AKC Rankings")

# Helper: write a literal text value into a cell without Excel coercing a
# date-shaped string into a real date serial, and without leaving a
# leftover number-format/quote-prefix style behind.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Drop the bold / centered / thin-border header formatting on row 1.
$ws.Range("A1:C1").ClearFormats()
$ws.Range("A1").Value = "Search term"
$ws.Range("B1").Value = "Results Ranking"
$ws.Range("C1").Value = "Date"

# Existing rows 2-3: same search term/result, date refreshed.
Set-TextValue $ws.Range("C2") "2025-11-04"
Set-TextValue $ws.Range("C3") "2025-11-04"

# New rows 4-11: latest run of the checker.
$newRows = @(
    @{ A = "Food Safety Course"; B = "Not Found"; Num = $false },
    @{ A = "Food Safety Course"; B = "Not Found"; Num = $false },
    @{ A = "Food Safety Course"; B = "Not Found"; Num = $false },
    @{ A = "Food Safety Course"; B = "Not Found"; Num = $false },
    @{ A = "Food Safety Course"; B = 7;           Num = $true  },
    @{ A = "Food Safety Course"; B = 7;           Num = $true  },
    @{ A = "Food Safety Course"; B = 7;           Num = $true  },
    @{ A = "BCSS Course";        B = 2;           Num = $true  }
)

$r = 4
foreach ($row in $newRows) {
    Set-TextValue $ws.Range("A$r") $row.A
    if ($row.Num) {
        $ws.Range("B$r").Value = $row.B
    } else {
        Set-TextValue $ws.Range("B$r") $row.B
    }
    Set-TextValue $ws.Range("C$r") "2025-11-04"
    $r++
}
